# "Final pre-selection show predictions"
#
# On the "Tournament Teams" sheet, the empty "Conference" placeholder
# column (B) that sat between the Conference (A) / Team (C) / AdjEM (D)
# columns for the "AUTO BID TEAMS" table (rows 1-33) is removed, shifting
# the Team values from column C into column B, and the AdjEM values from
# column D into column C (matching the layout already used by the
# "AT-LARGE TEAMS" table further down the sheet).
#
# The workbook is also repointed so that "Tournament Teams" becomes the
# selected/active sheet (instead of "Round of 32"), with the active cell
# at D7.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Tournament Teams")
$ws.Activate()

# Shift columns C (Team) and D (AdjEM) one column to the left, into B and
# C, for rows 1 through 33 - leaving column D empty afterwards.
for ($r = 1; $r -le 33; $r++) {
    $teamVal  = $ws.Cells.Item($r, 3).Value2
    $adjEmVal = $ws.Cells.Item($r, 4).Value2

    $ws.Cells.Item($r, 2).Value = $teamVal
    $ws.Cells.Item($r, 3).Value = $adjEmVal
    $ws.Cells.Item($r, 4).ClearContents()
}

# Make "Tournament Teams" the active sheet/tab, with D7 selected, and
# scrolled back to the top-left corner.
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("D7").Select()
